# Slide 3 ("Tipos básicos de datos y operadores"), shape "3 CuadroTexto":
# Restructure the bullet list - rename/split the first three items and add a
# new "Variables y operaciones lógicas" item (now last) plus a trailing blank
# bullet line.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)

# --- Paragraph 1: "Variables y operaciones lógicas" -> "Asignar valores a variables"
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)
$para1.Delete()
$tr = $sh.TextFrame.TextRange
$newFirst = $tr.Paragraphs(1, 1)
$newFirst.InsertBefore("Asignar valores a variables`r")

# --- Paragraph 2: "Variables y operaciones numéricas" -> split into "Variables " / "y operaciones numéricas"
$tr = $sh.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
$head2 = $para2.Characters(1, 10)
$head2.Text = "Variables "

# --- Paragraph 3: "Variables y operaciones con cadenas" -> split into "Variables y operaciones con " / "cadenas"
$tr = $sh.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
$tail3 = $para3.Characters(29, 7)
$tail3.Text = "cadenas"

# --- New paragraph 4 ("Variables y operaciones lógicas") + blank paragraph 5, both
#     inheriting paragraph 3's numbered-bullet style.
$tr = $sh.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
$para3.InsertAfter("`r`r")
$tr = $sh.TextFrame.TextRange
$para4 = $tr.Paragraphs(4, 1)
$para4.Text = "Variables y operaciones lógicas"

Write-Output $sh.TextFrame.TextRange.Text
